# Adds two new weekly price records for "Coliflor" (Vega Modelo de Temuco)
# by inserting two rows above the current row 180, shifting all the
# subsequent rows down by two (the sheet appears to be kept sorted by
# date, and these rows belong right after the existing 44329 entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 180; everything from old row 180
# onward shifts down to start at row 182.
$ws.Rows("180:181").Insert()

# New row 180
$ws.Range("A180").Value = 10
$ws.Range("B180").Value = "Vega Modelo de Temuco"
$ws.Range("C180").Value = "La Araucanía"
$ws.Range("D180").Value = 44510
$ws.Range("E180").Value = 9
$ws.Range("F180").Value = 100112008
$ws.Range("G180").Value = "Coliflor"
$ws.Range("H180").Value = "Sin especificar"
$ws.Range("I180").Value = "Primera"
$ws.Range("J180").Value = 600
$ws.Range("K180").Value = 900
$ws.Range("L180").Value = 900
$ws.Range("M180").Value = 900
$ws.Range("N180").Value = "$/unidad"
$ws.Range("O180").Value = "Región de O'Higgins"
$ws.Range("P180").Value = 900
$ws.Range("Q180").Value = 1
$ws.Range("R180").Value = "Hortaliza"

# New row 181
$ws.Range("A181").Value = 10
$ws.Range("B181").Value = "Vega Modelo de Temuco"
$ws.Range("C181").Value = "La Araucanía"
$ws.Range("D181").Value = 44510
$ws.Range("E181").Value = 9
$ws.Range("F181").Value = 100112008
$ws.Range("G181").Value = "Coliflor"
$ws.Range("H181").Value = "Sin especificar"
$ws.Range("I181").Value = "Primera"
$ws.Range("J181").Value = 500
$ws.Range("K181").Value = 800
$ws.Range("L181").Value = 800
$ws.Range("M181").Value = 800
$ws.Range("N181").Value = "$/unidad"
$ws.Range("O181").Value = "Región del Maule"
$ws.Range("P181").Value = 800
$ws.Range("Q181").Value = 1
$ws.Range("R181").Value = "Hortaliza"
